# Week 3 update: reposition/resize the Picture + caption TextBox on the
# "Comparison of data for model types - Poisson" slide, and tidy up the
# wording on the "Zero-inflated problems" slide.

$p = $ppt.ActivePresentation

# --- Slide 21: Picture 3 + TextBox 4 ---------------------------------
# Shape.Left/Top/Width/Height are Single-precision points (as in real
# PowerPoint COM), and the engine truncates when converting back to EMU,
# so the literals below are chosen to land exactly on the target EMU
# after the float32 round-trip (914400 EMU = 72 points = 1 inch).
$s21 = $p.Slides.Item(21)

$pic = $s21.Shapes.Item(2)
$pic.Left   = 276.63158480314956   # -> 3513221 EMU
$pic.Top    = 81.35897637795276    # -> 1033259 EMU
$pic.Width  = 363.0315748031496    # -> 4610501 EMU
$pic.Height = 463.3044094488189    # -> 5883966 EMU

$caption = $s21.Shapes.Item(3)
$caption.Left = 0.0                # -> 0 EMU
$caption.Top  = 510.9187501574803  # -> 6488668 EMU

# --- Slide 22: drop the "(how?)" aside -------------------------------
$s22 = $p.Slides.Item(22)
$body = $s22.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$fullText = $tr.Text

$oldSentence = " lizard data, examine the difference between a zero-inflated Poisson and a zero-inflated negative binomial. Compare these models (how?)"
$newSentence = " lizard data, examine the difference between a zero-inflated Poisson and a zero-inflated negative binomial. Compare these models."

$idx = $fullText.IndexOf($oldSentence)
$run = $tr.Characters($idx + 1, $oldSentence.Length)
$run.Text = $newSentence
